# Update scripts with new TPM values for the Agt-Mas1 LR-pair sheet.
# The sending/target cluster set grows from {FAPs, MuSCs} / {ECs, FAPs, MuSCs}
# to a full {ECs, FAPs, MuSCs} x {ECs, FAPs, MuSCs} cross-product (9 rows),
# and every NATMI-derived metric column is refreshed with the new TPM-based values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value2 = "ECs"
$ws.Cells.Item(2, 2).Value2 = "Agt"
$ws.Cells.Item(2, 3).Value2 = "Mas1"
$ws.Cells.Item(2, 4).Value2 = "ECs"
$ws.Cells.Item(2, 5).Value2 = 1
$ws.Cells.Item(2, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(2, 7).Value2 = 0.1206283333333333
$ws.Cells.Item(2, 8).Value2 = 0.361885
$ws.Cells.Item(2, 9).Value2 = 0.09993369694616584
$ws.Cells.Item(2, 10).Value2 = 0.09993369694616584
$ws.Cells.Item(2, 11).Value2 = 3
$ws.Cells.Item(2, 12).Value2 = 1
$ws.Cells.Item(2, 13).Value2 = 1.316251
$ws.Cells.Item(2, 14).Value2 = 3.948753
$ws.Cells.Item(2, 15).Value2 = 0.2348117009309627
$ws.Cells.Item(2, 16).Value2 = 0.2348117009309626
$ws.Cells.Item(2, 17).Value2 = 0.1587771643783333
$ws.Cells.Item(2, 18).Value2 = 1.428994479405
$ws.Cells.Item(2, 19).Value2 = 0.02346560136024855
$ws.Cells.Item(2, 20).Value2 = 0.02346560136024855

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value2 = "ECs"
$ws.Cells.Item(3, 2).Value2 = "Agt"
$ws.Cells.Item(3, 3).Value2 = "Mas1"
$ws.Cells.Item(3, 4).Value2 = "FAPs"
$ws.Cells.Item(3, 5).Value2 = 1
$ws.Cells.Item(3, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(3, 7).Value2 = 0.1206283333333333
$ws.Cells.Item(3, 8).Value2 = 0.361885
$ws.Cells.Item(3, 9).Value2 = 0.09993369694616584
$ws.Cells.Item(3, 10).Value2 = 0.09993369694616584
$ws.Cells.Item(3, 11).Value2 = 3
$ws.Cells.Item(3, 12).Value2 = 1
$ws.Cells.Item(3, 13).Value2 = 3.062384333333333
$ws.Cells.Item(3, 14).Value2 = 9.187152999999999
$ws.Cells.Item(3, 15).Value2 = 0.5463119680169907
$ws.Cells.Item(3, 16).Value2 = 0.5463119680169907
$ws.Cells.Item(3, 17).Value2 = 0.369410318156111
$ws.Cells.Item(3, 18).Value2 = 3.324692863405
$ws.Cells.Item(3, 19).Value2 = 0.05459497464987339
$ws.Cells.Item(3, 20).Value2 = 0.05459497464987339

# Row 4: ECs -> MuSCs
$ws.Cells.Item(4, 1).Value2 = "ECs"
$ws.Cells.Item(4, 2).Value2 = "Agt"
$ws.Cells.Item(4, 3).Value2 = "Mas1"
$ws.Cells.Item(4, 4).Value2 = "MuSCs"
$ws.Cells.Item(4, 5).Value2 = 1
$ws.Cells.Item(4, 6).Value2 = 0.3333333333333333
$ws.Cells.Item(4, 7).Value2 = 0.1206283333333333
$ws.Cells.Item(4, 8).Value2 = 0.361885
$ws.Cells.Item(4, 9).Value2 = 0.09993369694616584
$ws.Cells.Item(4, 10).Value2 = 0.09993369694616584
$ws.Cells.Item(4, 11).Value2 = 3
$ws.Cells.Item(4, 12).Value2 = 1
$ws.Cells.Item(4, 13).Value2 = 1.226924333333334
$ws.Cells.Item(4, 14).Value2 = 3.680773
$ws.Cells.Item(4, 15).Value2 = 0.2188763310520467
$ws.Cells.Item(4, 16).Value2 = 0.2188763310520466
$ws.Cells.Item(4, 17).Value2 = 0.1480018374561111
$ws.Cells.Item(4, 18).Value2 = 1.332016537105
$ws.Cells.Item(4, 19).Value2 = 0.0218731209360439
$ws.Cells.Item(4, 20).Value2 = 0.02187312093604389

# Row 5: FAPs -> ECs
$ws.Cells.Item(5, 1).Value2 = "FAPs"
$ws.Cells.Item(5, 2).Value2 = "Agt"
$ws.Cells.Item(5, 3).Value2 = "Mas1"
$ws.Cells.Item(5, 4).Value2 = "ECs"
$ws.Cells.Item(5, 5).Value2 = 2
$ws.Cells.Item(5, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(5, 7).Value2 = 0.4683593333333333
$ws.Cells.Item(5, 8).Value2 = 1.405078
$ws.Cells.Item(5, 9).Value2 = 0.3880090057275787
$ws.Cells.Item(5, 10).Value2 = 0.3880090057275787
$ws.Cells.Item(5, 11).Value2 = 3
$ws.Cells.Item(5, 12).Value2 = 1
$ws.Cells.Item(5, 13).Value2 = 1.316251
$ws.Cells.Item(5, 14).Value2 = 3.948753
$ws.Cells.Item(5, 15).Value2 = 0.2348117009309627
$ws.Cells.Item(5, 16).Value2 = 0.2348117009309626
$ws.Cells.Item(5, 17).Value2 = 0.6164784408593333
$ws.Cells.Item(5, 18).Value2 = 5.548305967734001
$ws.Cells.Item(5, 19).Value2 = 0.09110905461142439
$ws.Cells.Item(5, 20).Value2 = 0.09110905461142436

# Row 6: FAPs -> FAPs
$ws.Cells.Item(6, 1).Value2 = "FAPs"
$ws.Cells.Item(6, 2).Value2 = "Agt"
$ws.Cells.Item(6, 3).Value2 = "Mas1"
$ws.Cells.Item(6, 4).Value2 = "FAPs"
$ws.Cells.Item(6, 5).Value2 = 2
$ws.Cells.Item(6, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(6, 7).Value2 = 0.4683593333333333
$ws.Cells.Item(6, 8).Value2 = 1.405078
$ws.Cells.Item(6, 9).Value2 = 0.3880090057275787
$ws.Cells.Item(6, 10).Value2 = 0.3880090057275787
$ws.Cells.Item(6, 11).Value2 = 3
$ws.Cells.Item(6, 12).Value2 = 1
$ws.Cells.Item(6, 13).Value2 = 3.062384333333333
$ws.Cells.Item(6, 14).Value2 = 9.187152999999999
$ws.Cells.Item(6, 15).Value2 = 0.5463119680169907
$ws.Cells.Item(6, 16).Value2 = 0.5463119680169907
$ws.Cells.Item(6, 17).Value2 = 1.434296284770444
$ws.Cells.Item(6, 18).Value2 = 12.908666562934
$ws.Cells.Item(6, 19).Value2 = 0.2119739635273493
$ws.Cells.Item(6, 20).Value2 = 0.2119739635273493

# Row 7: FAPs -> MuSCs
$ws.Cells.Item(7, 1).Value2 = "FAPs"
$ws.Cells.Item(7, 2).Value2 = "Agt"
$ws.Cells.Item(7, 3).Value2 = "Mas1"
$ws.Cells.Item(7, 4).Value2 = "MuSCs"
$ws.Cells.Item(7, 5).Value2 = 2
$ws.Cells.Item(7, 6).Value2 = 0.6666666666666666
$ws.Cells.Item(7, 7).Value2 = 0.4683593333333333
$ws.Cells.Item(7, 8).Value2 = 1.405078
$ws.Cells.Item(7, 9).Value2 = 0.3880090057275787
$ws.Cells.Item(7, 10).Value2 = 0.3880090057275787
$ws.Cells.Item(7, 11).Value2 = 3
$ws.Cells.Item(7, 12).Value2 = 1
$ws.Cells.Item(7, 13).Value2 = 1.226924333333334
$ws.Cells.Item(7, 14).Value2 = 3.680773
$ws.Cells.Item(7, 15).Value2 = 0.2188763310520467
$ws.Cells.Item(7, 16).Value2 = 0.2188763310520466
$ws.Cells.Item(7, 17).Value2 = 0.5746414628104446
$ws.Cells.Item(7, 18).Value2 = 5.171773165294001
$ws.Cells.Item(7, 19).Value2 = 0.08492598758880497
$ws.Cells.Item(7, 20).Value2 = 0.08492598758880497

# Row 8: MuSCs -> ECs
$ws.Cells.Item(8, 1).Value2 = "MuSCs"
$ws.Cells.Item(8, 2).Value2 = "Agt"
$ws.Cells.Item(8, 3).Value2 = "Mas1"
$ws.Cells.Item(8, 4).Value2 = "ECs"
$ws.Cells.Item(8, 5).Value2 = 3
$ws.Cells.Item(8, 6).Value2 = 1
$ws.Cells.Item(8, 7).Value2 = 0.618096
$ws.Cells.Item(8, 8).Value2 = 1.854288
$ws.Cells.Item(8, 9).Value2 = 0.5120572973262555
$ws.Cells.Item(8, 10).Value2 = 0.5120572973262555
$ws.Cells.Item(8, 11).Value2 = 3
$ws.Cells.Item(8, 12).Value2 = 1
$ws.Cells.Item(8, 13).Value2 = 1.316251
$ws.Cells.Item(8, 14).Value2 = 3.948753
$ws.Cells.Item(8, 15).Value2 = 0.2348117009309627
$ws.Cells.Item(8, 16).Value2 = 0.2348117009309626
$ws.Cells.Item(8, 17).Value2 = 0.813569478096
$ws.Cells.Item(8, 18).Value2 = 7.322125302863999
$ws.Cells.Item(8, 19).Value2 = 0.1202370449592897
$ws.Cells.Item(8, 20).Value2 = 0.1202370449592897

# Row 9: MuSCs -> FAPs
$ws.Cells.Item(9, 1).Value2 = "MuSCs"
$ws.Cells.Item(9, 2).Value2 = "Agt"
$ws.Cells.Item(9, 3).Value2 = "Mas1"
$ws.Cells.Item(9, 4).Value2 = "FAPs"
$ws.Cells.Item(9, 5).Value2 = 3
$ws.Cells.Item(9, 6).Value2 = 1
$ws.Cells.Item(9, 7).Value2 = 0.618096
$ws.Cells.Item(9, 8).Value2 = 1.854288
$ws.Cells.Item(9, 9).Value2 = 0.5120572973262555
$ws.Cells.Item(9, 10).Value2 = 0.5120572973262555
$ws.Cells.Item(9, 11).Value2 = 3
$ws.Cells.Item(9, 12).Value2 = 1
$ws.Cells.Item(9, 13).Value2 = 3.062384333333333
$ws.Cells.Item(9, 14).Value2 = 9.187152999999999
$ws.Cells.Item(9, 15).Value2 = 0.5463119680169907
$ws.Cells.Item(9, 16).Value2 = 0.5463119680169907
$ws.Cells.Item(9, 17).Value2 = 1.892847506896
$ws.Cells.Item(9, 18).Value2 = 17.035627562064
$ws.Cells.Item(9, 19).Value2 = 0.279743029839768
$ws.Cells.Item(9, 20).Value2 = 0.279743029839768

# Row 10: MuSCs -> MuSCs
$ws.Cells.Item(10, 1).Value2 = "MuSCs"
$ws.Cells.Item(10, 2).Value2 = "Agt"
$ws.Cells.Item(10, 3).Value2 = "Mas1"
$ws.Cells.Item(10, 4).Value2 = "MuSCs"
$ws.Cells.Item(10, 5).Value2 = 3
$ws.Cells.Item(10, 6).Value2 = 1
$ws.Cells.Item(10, 7).Value2 = 0.618096
$ws.Cells.Item(10, 8).Value2 = 1.854288
$ws.Cells.Item(10, 9).Value2 = 0.5120572973262555
$ws.Cells.Item(10, 10).Value2 = 0.5120572973262555
$ws.Cells.Item(10, 11).Value2 = 3
$ws.Cells.Item(10, 12).Value2 = 1
$ws.Cells.Item(10, 13).Value2 = 1.226924333333334
$ws.Cells.Item(10, 14).Value2 = 3.680773
$ws.Cells.Item(10, 15).Value2 = 0.2188763310520467
$ws.Cells.Item(10, 16).Value2 = 0.2188763310520466
$ws.Cells.Item(10, 17).Value2 = 0.758357022736
$ws.Cells.Item(10, 18).Value2 = 6.825213204624
$ws.Cells.Item(10, 19).Value2 = 0.1120772225271978
$ws.Cells.Item(10, 20).Value2 = 0.1120772225271978
